# Bring "Sheet2 - Numbers" in line with the reader/sheet2 fixture:
#  - drop the stray G5 string cell that didn't belong on this sheet
#  - add a new AA column (AA1:AA30 = 100..129)
#  - select the new AA column and make Sheet2 the active/visible tab
#  - Sheet4's page setup picks up an explicit paper size

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2 - Numbers")

# Remove the erroneous G5 entry (leftover from Sheet1) that doesn't belong here.
$ws2.Range("G5").ClearContents() | Out-Null

# Populate the new AA column with values 100-129 for rows 1-30.
for ($row = 1; $row -le 30; $row++) {
    $ws2.Cells.Item($row, 27).Value = 99 + $row
}

# Make Sheet2 the active sheet/tab and select the freshly written column.
$ws2.Activate()
$ws2.Range("AA1:AA30").Select() | Out-Null

# Sheet4 now reports an explicit (Letter) paper size instead of "unset".
$ws4 = $wb.Worksheets.Item("Sheet4 - Dates")
$ws4.PageSetup.PaperSize = 9
